$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (week of 2021-12-10) is inserted at row 30, pushing the
# existing rows 30-81 down by one (old row 81 becomes row 82).
$ws.Rows.Item(30).Insert()

# Fill in the new row 30 with the inserted record's data. Columns that stay
# constant across this market/product (A, B, C, E, F, G, H, I, R) are copied
# from the row immediately below (the former row 30, now row 31).
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Cells.Item(30, 4).Value = 44540
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = 100112052
$ws.Cells.Item(30, 7).Value = "Albahaca"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 130
$ws.Cells.Item(30, 11).Value = 7000
$ws.Cells.Item(30, 12).Value = 7000
$ws.Cells.Item(30, 13).Value = 7000
$ws.Cells.Item(30, 14).Value = "`$/docena de matas"
$ws.Cells.Item(30, 15).Value = "Región Metropolitana"
$ws.Cells.Item(30, 16).Value = 1167
$ws.Cells.Item(30, 17).Value = 6
$ws.Cells.Item(30, 18).Value = "Hortaliza"
